$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 47.875  # H5: was 41.166668
$ws.Cells.Item(5, 9).Value = 47  # I5: was 35.2
$ws.Cells.Item(5, 10).Value = 48.166668  # J5: was 45.42857
$ws.Cells.Item(5, 11).Value = 47  # K5: was 35.2
$ws.Cells.Item(5, 12).Value = 48.166668  # L5: was 45.42857
$ws.Cells.Item(5, 13).Value = 68  # M5: was 79.8
$ws.Cells.Item(5, 14).Value = -278.166668  # N5: was -275.42857
$ws.Cells.Item(18, 8).Value = 836.4  # H18: was 927.6667
$ws.Cells.Item(18, 9).Value = 836.4  # I18: was 927.6667
$ws.Cells.Item(18, 11).Value = 836.4  # K18: was 927.6667
$ws.Cells.Item(18, 13).Value = -552.4  # M18: was -643.6667
$ws.Cells.Item(19, 8).Value = 1206.8572  # H19: was 1118.4348
$ws.Cells.Item(19, 10).Value = 1259  # J19: was 1140.2222
$ws.Cells.Item(19, 12).Value = 1259  # L19: was 1140.2222
$ws.Cells.Item(19, 14).Value = -1609  # N19: was -1490.2222
$ws.Cells.Item(113, 8).Value = 2250.75  # H113: was 2200.4
$ws.Cells.Item(113, 9).Value = 2250.75  # I113: was 2200.4
$ws.Cells.Item(113, 11).Value = 2250.75  # K113: was 2200.4
$ws.Cells.Item(113, 13).Value = 1003.25  # M113: was 1053.6
$ws.Cells.Item(132, 8).Value = 1206.3334  # H132: was 1221.4849
$ws.Cells.Item(132, 9).Value = 1214.5077  # I132: was 1229.8923
$ws.Cells.Item(132, 11).Value = 3643.5231  # K132: was 3689.6769
$ws.Cells.Item(132, 13).Value = -1113.5231  # M132: was -1159.6769
$ws.Cells.Item(137, 8).Value = 393562.5  # H137: was 364060.44
$ws.Cells.Item(137, 9).Value = 1240.069  # I137: was 1173.0646
$ws.Cells.Item(137, 10).Value = 1815731.2  # J137: was 1614005.9
$ws.Cells.Item(137, 11).Value = 3720.207  # K137: was 3519.1938
$ws.Cells.Item(137, 12).Value = 5447193.6  # L137: was 4842017.699999999
$ws.Cells.Item(137, 13).Value = -1170.207  # M137: was -969.1938
$ws.Cells.Item(137, 14).Value = -5452293.6  # N137: was -4847117.699999999
$ws.Cells.Item(138, 8).Value = 47669356  # H138: was 47669324
$ws.Cells.Item(138, 10).Value = 76926330  # J138: was 76926280
$ws.Cells.Item(138, 12).Value = 230778990  # L138: was 230778840
$ws.Cells.Item(138, 14).Value = -230789270  # N138: was -230789120

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(46, 8).Value = 26632  # H46: was 24998
$ws.Cells.Item(46, 10).Value = 29948  # J46: was 29996
$ws.Cells.Item(46, 12).Value = 29948  # L46: was 29996
$ws.Cells.Item(46, 14).Value = -30586  # N46: was -30634
$ws.Cells.Item(63, 8).Value = 4497.6665  # H63: was 4633.8184
$ws.Cells.Item(63, 9).Value = 3795.4  # I63: was 3994.25
$ws.Cells.Item(63, 11).Value = 3795.4  # K63: was 3994.25
$ws.Cells.Item(63, 13).Value = -3109.4  # M63: was -3308.25
$ws.Cells.Item(66, 8).Value = 4497.6665  # H66: was 4633.8184
$ws.Cells.Item(66, 9).Value = 3795.4  # I66: was 3994.25
$ws.Cells.Item(66, 11).Value = 18977  # K66: was 19971.25
$ws.Cells.Item(66, 13).Value = -15545  # M66: was -16539.25
$ws.Cells.Item(97, 8).Value = 617.25  # H97: was 534.35
$ws.Cells.Item(97, 9).Value = 648.8  # I97: was 560
$ws.Cells.Item(97, 10).Value = 459.5  # J97: was 431.75
$ws.Cells.Item(97, 11).Value = 648.8  # K97: was 560
$ws.Cells.Item(97, 12).Value = 459.5  # L97: was 431.75
$ws.Cells.Item(97, 13).Value = -152.8  # M97: was -64
$ws.Cells.Item(97, 14).Value = -1451.5  # N97: was -1423.75
$ws.Cells.Item(110, 8).Value = 1553.3334  # H110: was 1640
$ws.Cells.Item(110, 9).Value = 1071.1111  # I110: was 1130
$ws.Cells.Item(110, 11).Value = 1071.1111  # K110: was 1130
$ws.Cells.Item(110, 13).Value = 973.8888999999999  # M110: was 915
$ws.Cells.Item(139, 8).Value = 155306.28  # H139: was 164524
$ws.Cells.Item(139, 10).Value = 155306.28  # J139: was 164524
$ws.Cells.Item(139, 12).Value = 155306.28  # L139: was 164524
$ws.Cells.Item(139, 14).Value = -165586.28  # N139: was -174804

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 5720.5  # H86: was 6086
$ws.Cells.Item(86, 9).Value = 3242  # I86: was 4290.4
$ws.Cells.Item(86, 10).Value = 6649.9375  # J86: was 6647.125
$ws.Cells.Item(86, 11).Value = 3242  # K86: was 4290.4
$ws.Cells.Item(86, 12).Value = 6649.9375  # L86: was 6647.125
$ws.Cells.Item(86, 13).Value = -2119  # M86: was -3167.4
$ws.Cells.Item(86, 14).Value = -8895.9375  # N86: was -8893.125
$ws.Cells.Item(89, 8).Value = 5720.5  # H89: was 6086
$ws.Cells.Item(89, 9).Value = 3242  # I89: was 4290.4
$ws.Cells.Item(89, 10).Value = 6649.9375  # J89: was 6647.125
$ws.Cells.Item(89, 11).Value = 16210  # K89: was 21452
$ws.Cells.Item(89, 12).Value = 33249.6875  # L89: was 33235.625
$ws.Cells.Item(89, 13).Value = -10594  # M89: was -15836
$ws.Cells.Item(89, 14).Value = -44481.6875  # N89: was -44467.625
$ws.Cells.Item(105, 8).Value = 52922.15  # H105: was 55770.684
$ws.Cells.Item(105, 9).Value = 61690.766  # I105: was 69862.87
$ws.Cells.Item(105, 10).Value = 3233.3333  # J105: was 2925
$ws.Cells.Item(105, 11).Value = 61690.766  # K105: was 69862.87
$ws.Cells.Item(105, 12).Value = 3233.3333  # L105: was 2925
$ws.Cells.Item(105, 13).Value = -59943.766  # M105: was -68115.87
$ws.Cells.Item(105, 14).Value = -6727.3333  # N105: was -6419
$ws.Cells.Item(107, 8).Value = 1859.6666  # H107: was 1929.0435
$ws.Cells.Item(107, 9).Value = 1715.2941  # I107: was 1806
$ws.Cells.Item(107, 11).Value = 1715.2941  # K107: was 1806
$ws.Cells.Item(107, 13).Value = 204.7058999999999  # M107: was 114

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 53000  # H20: was 59500
$ws.Cells.Item(20, 10).Value = 53000  # J20: was 59500
$ws.Cells.Item(20, 12).Value = 53000  # L20: was 59500
$ws.Cells.Item(20, 14).Value = -53472  # N20: was -59972
$ws.Cells.Item(25, 8).Value = 421.6  # H25: was 569.3333
$ws.Cells.Item(25, 9).Value = 369.33334  # I25: was 708
$ws.Cells.Item(25, 11).Value = 369.33334  # K25: was 708
$ws.Cells.Item(25, 13).Value = -195.33334  # M25: was -534
$ws.Cells.Item(30, 8).Value = 53000  # H30: was 59500
$ws.Cells.Item(30, 10).Value = 53000  # J30: was 59500
$ws.Cells.Item(30, 12).Value = 53000  # L30: was 59500
$ws.Cells.Item(30, 14).Value = -53182  # N30: was -59682
$ws.Cells.Item(59, 8).Value = 75000  # H59: was 169998
$ws.Cells.Item(59, 9).Value = 40000  # I59: was 0
$ws.Cells.Item(59, 10).Value = 110000  # J59: was 169998
$ws.Cells.Item(59, 11).Value = 40000  # K59: was 0
$ws.Cells.Item(59, 12).Value = 110000  # L59: was 169998
$ws.Cells.Item(59, 13).Value = -38855  # M59: was None
$ws.Cells.Item(59, 14).Value = -112290  # N59: was -172288
$ws.Cells.Item(128, 8).Value = 53000  # H128: was 59500
$ws.Cells.Item(128, 10).Value = 53000  # J128: was 59500
$ws.Cells.Item(128, 12).Value = 53000  # L128: was 59500
$ws.Cells.Item(128, 14).Value = -62960  # N128: was -69460
$ws.Cells.Item(132, 8).Value = 2927896.8  # H132: was 3252896.5
$ws.Cells.Item(132, 9).Value = 3252976.8  # I132: was 3794656.5
$ws.Cells.Item(132, 11).Value = 9758930.399999999  # K132: was 11383969.5
$ws.Cells.Item(132, 13).Value = -9756400.399999999  # M132: was -11381439.5
$ws.Cells.Item(134, 8).Value = 2014102.1  # H134: was 2071628.6
$ws.Cells.Item(134, 9).Value = 2553077.8  # I134: was 2647611.2
$ws.Cells.Item(134, 11).Value = 7659233.399999999  # K134: was 7942833.600000001
$ws.Cells.Item(134, 13).Value = -7656698.399999999  # M134: was -7940298.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 74999.75  # H37: was 74999.60000000001
$ws.Cells.Item(37, 10).Value = 74999.75  # J37: was 74999.60000000001
$ws.Cells.Item(37, 12).Value = 224999.25  # L37: was 224998.8
$ws.Cells.Item(37, 14).Value = -225223.25  # N37: was -225222.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 2418.9524  # H21: was 2179.0833
$ws.Cells.Item(21, 10).Value = 3159.6  # J21: was 2162.25
$ws.Cells.Item(21, 12).Value = 3159.6  # L21: was 2162.25
$ws.Cells.Item(21, 14).Value = -3505.6  # N21: was -2508.25
$ws.Cells.Item(30, 8).Value = 2418.9524  # H30: was 2179.0833
$ws.Cells.Item(30, 10).Value = 3159.6  # J30: was 2162.25
$ws.Cells.Item(30, 12).Value = 3159.6  # L30: was 2162.25
$ws.Cells.Item(30, 14).Value = -3369.6  # N30: was -2372.25
$ws.Cells.Item(47, 8).Value = 6997  # H47: was 6170.75
$ws.Cells.Item(47, 10).Value = 6997  # J47: was 6170.75
$ws.Cells.Item(47, 12).Value = 6997  # L47: was 6170.75
$ws.Cells.Item(47, 14).Value = -8133  # N47: was -7306.75
$ws.Cells.Item(113, 8).Value = 4439.4546  # H113: was 4633.4
$ws.Cells.Item(113, 10).Value = 4585.3335  # J113: was 5002.4
$ws.Cells.Item(113, 12).Value = 4585.3335  # L113: was 5002.4
$ws.Cells.Item(113, 14).Value = -8925.333500000001  # N113: was -9342.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 10554.538  # H7: was 28502.25
$ws.Cells.Item(7, 9).Value = 36003  # I7: was 53004.5
$ws.Cells.Item(7, 10).Value = 2920  # J7: was 4000
$ws.Cells.Item(7, 11).Value = 36003  # K7: was 53004.5
$ws.Cells.Item(7, 12).Value = 2920  # L7: was 4000
$ws.Cells.Item(7, 13).Value = -35891  # M7: was -52892.5
$ws.Cells.Item(7, 14).Value = -3144  # N7: was -4224
$ws.Cells.Item(16, 8).Value = 1544.8  # H16: was 1423.2778
$ws.Cells.Item(16, 9).Value = 1247.6666  # I16: was 1289.6923
$ws.Cells.Item(16, 10).Value = 2733.3333  # J16: was 1770.6
$ws.Cells.Item(16, 11).Value = 1247.6666  # K16: was 1289.6923
$ws.Cells.Item(16, 12).Value = 2733.3333  # L16: was 1770.6
$ws.Cells.Item(16, 13).Value = -1077.6666  # M16: was -1119.6923
$ws.Cells.Item(16, 14).Value = -3073.3333  # N16: was -2110.6
$ws.Cells.Item(29, 8).Value = 0  # H29: was 5500
$ws.Cells.Item(29, 10).Value = 0  # J29: was 5500
$ws.Cells.Item(29, 12).Value = 0  # L29: was 5500
$ws.Cells.Item(29, 14).ClearContents()  # N29: was -6090
$ws.Cells.Item(40, 8).Value = 6288275.5  # H40: was 1718176.8
$ws.Cells.Item(40, 9).Value = 147416.86  # I40: was 40029.215
$ws.Cells.Item(40, 10).Value = 27781280  # J40: was 11115803
$ws.Cells.Item(40, 11).Value = 147416.86  # K40: was 40029.215
$ws.Cells.Item(40, 12).Value = 27781280  # L40: was 11115803
$ws.Cells.Item(40, 13).Value = -147280.86  # M40: was -39893.215
$ws.Cells.Item(40, 14).Value = -27781552  # N40: was -11116075
$ws.Cells.Item(61, 8).Value = 7501.6  # H61: was 7627
$ws.Cells.Item(61, 9).Value = 7501.6  # I61: was 7627
$ws.Cells.Item(61, 11).Value = 7501.6  # K61: was 7627
$ws.Cells.Item(61, 13).Value = -7299.6  # M61: was -7425
$ws.Cells.Item(111, 8).Value = 81000  # H111: was 83000
$ws.Cells.Item(111, 10).Value = 81000  # J111: was 83000
$ws.Cells.Item(111, 12).Value = 81000  # L111: was 83000
$ws.Cells.Item(111, 14).Value = -89180  # N111: was -91180
$ws.Cells.Item(113, 8).Value = 7501.6  # H113: was 7627
$ws.Cells.Item(113, 9).Value = 7501.6  # I113: was 7627
$ws.Cells.Item(113, 11).Value = 7501.6  # K113: was 7627
$ws.Cells.Item(113, 13).Value = -5331.6  # M113: was -5457
$ws.Cells.Item(116, 8).Value = 274750  # H116: was 266333.34
$ws.Cells.Item(116, 10).Value = 274750  # J116: was 266333.34
$ws.Cells.Item(116, 12).Value = 274750  # L116: was 266333.34
$ws.Cells.Item(116, 14).Value = -283928  # N116: was -275511.34
$ws.Cells.Item(126, 8).Value = 10554.538  # H126: was 28502.25
$ws.Cells.Item(126, 9).Value = 36003  # I126: was 53004.5
$ws.Cells.Item(126, 10).Value = 2920  # J126: was 4000
$ws.Cells.Item(126, 11).Value = 108009  # K126: was 159013.5
$ws.Cells.Item(126, 12).Value = 8760  # L126: was 12000
$ws.Cells.Item(126, 13).Value = -105539  # M126: was -156543.5
$ws.Cells.Item(126, 14).Value = -13700  # N126: was -16940
$ws.Cells.Item(132, 8).Value = 3431.7273  # H132: was 3983.2222
$ws.Cells.Item(132, 9).Value = 3194.3333  # I132: was 3835.5715
$ws.Cells.Item(132, 11).Value = 9582.999899999999  # K132: was 11506.7145
$ws.Cells.Item(132, 13).Value = -7052.999899999999  # M132: was -8976.7145
$ws.Cells.Item(134, 8).Value = 149999  # H134: was 0
$ws.Cells.Item(134, 10).Value = 149999  # J134: was 0
$ws.Cells.Item(134, 12).Value = 149999  # L134: was 0
$ws.Cells.Item(134, 14).Value = -160139  # N134: was None
$ws.Cells.Item(135, 8).Value = 100000  # H135: was 119999.5
$ws.Cells.Item(135, 10).Value = 100000  # J135: was 119999.5
$ws.Cells.Item(135, 12).Value = 100000  # L135: was 119999.5
$ws.Cells.Item(135, 14).Value = -110140  # N135: was -130139.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 0  # H32: was 1000
$ws.Cells.Item(32, 10).Value = 0  # J32: was 1000
$ws.Cells.Item(32, 12).Value = 0  # L32: was 1000
$ws.Cells.Item(32, 14).ClearContents()  # N32: was -1634
$ws.Cells.Item(74, 8).Value = 36124.25  # H74: was 37377.5
$ws.Cells.Item(74, 10).Value = 36124.25  # J74: was 37377.5
$ws.Cells.Item(74, 12).Value = 36124.25  # L74: was 37377.5
$ws.Cells.Item(74, 14).Value = -37996.25  # N74: was -39249.5
$ws.Cells.Item(77, 8).Value = 36124.25  # H77: was 37377.5
$ws.Cells.Item(77, 10).Value = 36124.25  # J77: was 37377.5
$ws.Cells.Item(77, 12).Value = 108372.75  # L77: was 112132.5
$ws.Cells.Item(77, 14).Value = -117732.75  # N77: was -121492.5
$ws.Cells.Item(96, 8).Value = 8792180  # H96: was 4397169
$ws.Cells.Item(96, 9).Value = 38667.668  # I96: was 20333.666
$ws.Cells.Item(96, 10).Value = 17545692  # J96: was 8774005
$ws.Cells.Item(96, 11).Value = 38667.668  # K96: was 20333.666
$ws.Cells.Item(96, 12).Value = 17545692  # L96: was 8774005
$ws.Cells.Item(96, 13).Value = -37294.668  # M96: was -18960.666
$ws.Cells.Item(96, 14).Value = -17548438  # N96: was -8776751
$ws.Cells.Item(108, 8).Value = 0  # H108: was 100000
$ws.Cells.Item(108, 10).Value = 0  # J108: was 100000
$ws.Cells.Item(108, 12).Value = 0  # L108: was 100000
$ws.Cells.Item(108, 14).ClearContents()  # N108: was -107680
$ws.Cells.Item(109, 8).Value = 0  # H109: was 49999
$ws.Cells.Item(109, 10).Value = 0  # J109: was 49999
$ws.Cells.Item(109, 12).Value = 0  # L109: was 49999
$ws.Cells.Item(109, 14).ClearContents()  # N109: was -52773
$ws.Cells.Item(122, 8).Value = 2948.5833  # H122: was 3020.261
$ws.Cells.Item(122, 9).Value = 2906  # I122: was 2986.3
$ws.Cells.Item(122, 11).Value = 8718  # K122: was 8958.900000000001
$ws.Cells.Item(122, 13).Value = -6268  # M122: was -6508.900000000001
$ws.Cells.Item(126, 8).Value = 4206.5454  # H126: was 2987.2778
$ws.Cells.Item(126, 9).Value = 3263.8572  # I126: was 2705.111
$ws.Cells.Item(126, 10).Value = 5856.25  # J126: was 3269.4443
$ws.Cells.Item(126, 11).Value = 9791.571599999999  # K126: was 8115.333
$ws.Cells.Item(126, 12).Value = 17568.75  # L126: was 9808.332900000001
$ws.Cells.Item(126, 13).Value = -7321.571599999999  # M126: was -5645.333
$ws.Cells.Item(126, 14).Value = -22508.75  # N126: was -14748.3329
